{"js": "const pairs = [\n  [\"2025-07-15 Tuesday\", \"2025-07-16 Wednesday\"],\n  [\"497\u00f74=124, 1\", \"953\u00f75=190, 3\"],\n  [\"839\u00f78=104, 7\", \"849\u00f79=94, 3\"],\n  [\"246\u00f73=82, 0\", \"641\u00f79=71, 2\"],\n  [\"663\u00f75=132, 3\", \"458\u00f77=65, 3\"],\n  [\"332\u00f74=83, 0\", \"582\u00f74=145, 2\"],\n  [\"820\u00f72=410, 0\", \"110\u00f76=18, 2\"],\n  [\"597\u00f78=74, 5\", \"214\u00f78=26, 6\"],\n  [\"983\u00f73=327, 2\", \"611\u00f73=203, 2\"],\n  [\"443\u00f75=88, 3\", \"751\u00f74=187, 3\"],\n  [\"409\u00f78=51, 1\", \"180\u00f76=30, 0\"],\n  [\"255\u00f72=127, 1\", \"618\u00f76=103, 0\"],\n  [\"851\u00f73=283, 2\", \"152\u00f75=30, 2\"],\n  [\"223\u00f75=44, 3\", \"340\u00f72=170, 0\"],\n  [\"209\u00f72=104, 1\", \"356\u00f76=59, 2\"],\n  [\"455\u00f73=151, 2\", \"185\u00f74=46, 1\"],\n  [\"301\u00f72=150, 1\", \"142\u00f73=47, 1\"],\n  [\"517\u00f74=129, 1\", \"719\u00f79=79, 8\"],\n  [\"741\u00f77=105, 6\", \"901\u00f72=450, 1\"],\n  [\"387\u00f76=64, 3\", \"123\u00f77=17, 4\"],\n  [\"319\u00f77=45, 4\", \"988\u00f73=329, 1\"],\n  [\"326\u00f75=65, 1\", \"416\u00f76=69, 2\"],\n  [\"134\u00f78=16, 6\", \"904\u00f77=129, 1\"],\n  [\"305\u00f78=38, 1\", \"163\u00f72=81, 1\"],\n  [\"193\u00f76=32, 1\", \"512\u00f77=73, 1\"],\n  [\"481\u00f74=120, 1\", \"219\u00f79=24, 3\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Not found: \" + oldText);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2025-07-15 Tuesday\", \"2025-07-16 Wednesday\"),\n    @(\"497\u00f74=124, 1\", \"953\u00f75=190, 3\"),\n    @(\"839\u00f78=104, 7\", \"849\u00f79=94, 3\"),\n    @(\"246\u00f73=82, 0\", \"641\u00f79=71, 2\"),\n    @(\"663\u00f75=132, 3\", \"458\u00f77=65, 3\"),\n    @(\"332\u00f74=83, 0\", \"582\u00f74=145, 2\"),\n    @(\"820\u00f72=410, 0\", \"110\u00f76=18, 2\"),\n    @(\"597\u00f78=74, 5\", \"214\u00f78=26, 6\"),\n    @(\"983\u00f73=327, 2\", \"611\u00f73=203, 2\"),\n    @(\"443\u00f75=88, 3\", \"751\u00f74=187, 3\"),\n    @(\"409\u00f78=51, 1\", \"180\u00f76=30, 0\"),\n    @(\"255\u00f72=127, 1\", \"618\u00f76=103, 0\"),\n    @(\"851\u00f73=283, 2\", \"152\u00f75=30, 2\"),\n    @(\"223\u00f75=44, 3\", \"340\u00f72=170, 0\"),\n    @(\"209\u00f72=104, 1\", \"356\u00f76=59, 2\"),\n    @(\"455\u00f73=151, 2\", \"185\u00f74=46, 1\"),\n    @(\"301\u00f72=150, 1\", \"142\u00f73=47, 1\"),\n    @(\"517\u00f74=129, 1\", \"719\u00f79=79, 8\"),\n    @(\"741\u00f77=105, 6\", \"901\u00f72=450, 1\"),\n    @(\"387\u00f76=64, 3\", \"123\u00f77=17, 4\"),\n    @(\"319\u00f77=45, 4\", \"988\u00f73=329, 1\"),\n    @(\"326\u00f75=65, 1\", \"416\u00f76=69, 2\"),\n    @(\"134\u00f78=16, 6\", \"904\u00f77=129, 1\"),\n    @(\"305\u00f78=38, 1\", \"163\u00f72=81, 1\"),\n    @(\"193\u00f76=32, 1\", \"512\u00f77=73, 1\"),\n    @(\"481\u00f74=120, 1\", \"219\u00f79=24, 3\"),\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $found = $range.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n    if (-not $found) {\n        Write-Output \"NOT FOUND: $old\"\n    }\n}\n"}
